$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers but must remain stored as
# text (matching the original inlineStr cell type). We temporarily force a
# text number format so Excel does not auto-convert the assigned string into
# a numeric value (which would also strip meaningful trailing zeros), then
# restore the cell style back to Normal so no stray formatting is left behind.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "43.422.36"
$ws.Range("E2").Value = "  +2.82%  "
$ws.Range("D3").Value = "2.311.66"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "311.58"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("D6").Value = "102.30"
$ws.Range("E6").Value = "  +5.31%  "
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  +7.27%  "
$ws.Range("D10").Value = "35.81"
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("D11").Value = "0.0815"
$ws.Range("E11").Value = "  +2.91%  "
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "2.669.45"
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("D15").Value = "14.98"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").Value = "2.314.07"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").Value = "0.811"
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("D18").Value = "43.332.79"
$ws.Range("E18").Value = "  +2.91%  "
$ws.Range("D19").Value = "12.38"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "0.0₃0934"
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("D21").Value = "6.18"
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "241.70"
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("D24").Value = "2.62"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("D28").Value = "24.64"
$ws.Range("E28").Value = "  +4.50%  "
$ws.Range("D29").Value = "36.81"
$ws.Range("E29").Value = "  -3.04%  "
$ws.Range("D30").Value = "9.64"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "167.74"
$ws.Range("E32").Value = "  +3.29%  "
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "0.0745"
$ws.Range("E35").Value = "  +0.79%  "
$ws.Range("E36").Value = "  +5.53%  "
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("E38").Value = "  -2.82%  "
$ws.Range("E39").Value = "  +4.16%  "
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("E41").Value = "  +1.55%  "
$ws.Range("E42").Value = "  +7.03%  "
$ws.Range("D43").Value = "2.32"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("E44").Value = "  +2.72%  "
$ws.Range("D45").Value = "1.974.39"
$ws.Range("D46").Value = "19.25"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("D48").Value = "9.91"
$ws.Range("E48").Value = "  +0.59%  "
$ws.Range("D49").Value = "55.76"
$ws.Range("E49").Value = "  +3.64%  "
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("E51").Value = "  +7.01%  "

# Restore default (Normal) style on the cells we temporarily formatted as text
# so the saved workbook does not retain an unused custom number format.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"

